# Applies the latest cryptos.xlsx scrape: refreshed Price (D) and
# Volume(1h) (E) figures, and a reordering of the RenderToken /
# InternetComputer(DFINITY) rows (30-31), swapping their Coin/Link/
# Price/Volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as plain text, even when the text looks like a
# number (e.g. "58.266.98", "529.30"). A leading apostrophe forces
# Excel to store it as text instead of auto-converting it to a
# number; resetting the style back to "Normal" afterwards drops the
# quote-prefix marker so the cell keeps its original (unstyled) look.
function Set-TextValue($address, $text) {
    $ws.Range($address).Value = "'" + $text
    $ws.Range($address).Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '58.266.98'
$ws.Range('E2').Value = '  -0.95%  '
# Row 3
Set-TextValue 'D3' '3.123.51'
$ws.Range('E3').Value = '  +1.05%  '
# Row 4
$ws.Range('E4').Value = '  -0.06%  '
# Row 5
Set-TextValue 'D5' '529.30'
$ws.Range('E5').Value = '  +1.47%  '
# Row 6
Set-TextValue 'D6' '142.63'
$ws.Range('E6').Value = '  -0.80%  '
# Row 7
$ws.Range('E7').Value = '  +0.06%  '
# Row 8
Set-TextValue 'D8' '3.121.44'
$ws.Range('E8').Value = '  +1.10%  '
# Row 9
$ws.Range('E9').Value = '  +1.38%  '
# Row 10
Set-TextValue 'D10' '7.18'
$ws.Range('E10').Value = '  -2.49%  '
# Row 11
$ws.Range('E11').Value = '  -0.19%  '
# Row 12
$ws.Range('E12').Value = '  +2.84%  '
# Row 13
Set-TextValue 'D13' '3.654.25'
$ws.Range('E13').Value = '  +0.84%  '
# Row 14
$ws.Range('E14').Value = '  +3.43%  '
# Row 15
Set-TextValue 'D15' '25.85'
$ws.Range('E15').Value = '  -3.38%  '
# Row 16
Set-TextValue 'D16' '0.0000166'
$ws.Range('E16').Value = '  -0.23%  '
# Row 17
Set-TextValue 'D17' '58.299.29'
$ws.Range('E17').Value = '  -0.91%  '
# Row 18
Set-TextValue 'D18' '3.129.74'
$ws.Range('E18').Value = '  +0.85%  '
# Row 19
Set-TextValue 'D19' '6.13'
$ws.Range('E19').Value = '  -0.18%  '
# Row 20
Set-TextValue 'D20' '12.83'
$ws.Range('E20').Value = '  -0.54%  '
# Row 21
Set-TextValue 'D21' '8.01'
$ws.Range('E21').Value = '  -1.12%  '
# Row 22
Set-TextValue 'D22' '343.91'
$ws.Range('E22').Value = '  +0.02%  '
# Row 23
$ws.Range('E23').Value = '  +0.11%  '
# Row 24
$ws.Range('E24').Value = '  +2.10%  '
# Row 25
Set-TextValue 'D25' '67.91'
$ws.Range('E25').Value = '  +3.31%  '
# Row 26
$ws.Range('E26').Value = '  -0.91%  '
# Row 27
$ws.Range('E27').Value = '  +0.02%  '
# Row 28
Set-TextValue 'D28' '0.0₃0932'
$ws.Range('E28').Value = '  +1.21%  '
# Row 29
$ws.Range('E29').Value = '  +0.12%  '
# Row 30
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D30' '7.33'
$ws.Range('E30').Value = '  +1.13%  '
# Row 31
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D31' '6.39'
$ws.Range('E31').Value = '  -3.21%  '
# Row 32
$ws.Range('E32').Value = '  +2.40%  '
# Row 33
Set-TextValue 'D33' '21.10'
$ws.Range('E33').Value = '  +0.62%  '
# Row 34
Set-TextValue 'D34' '1.20'
$ws.Range('E34').Value = '  -0.55%  '
# Row 35
Set-TextValue 'D35' '158.20'
$ws.Range('E35').Value = '  +2.23%  '
# Row 36
Set-TextValue 'D36' '4.66'
$ws.Range('E36').Value = '  +0.83%  '
# Row 37
Set-TextValue 'D37' '6.22'
$ws.Range('E37').Value = '  +1.44%  '
# Row 38
Set-TextValue 'D38' '26.44'
$ws.Range('E38').Value = '  -1.35%  '
# Row 39
Set-TextValue 'D39' '1.26'
$ws.Range('E39').Value = '  -3.16%  '
# Row 40
Set-TextValue 'D40' '1.65'
$ws.Range('E40').Value = '  +13.44%  '
# Row 41
Set-TextValue 'D41' '0.0671'
$ws.Range('E41').Value = '  -2.17%  '
# Row 42
Set-TextValue 'D42' '4.02'
$ws.Range('E42').Value = '  +2.39%  '
# Row 43
Set-TextValue 'D43' '0.693'
$ws.Range('E43').Value = '  +4.27%  '
# Row 44
Set-TextValue 'D44' '3.160.54'
$ws.Range('E44').Value = '  +0.87%  '
# Row 45
Set-TextValue 'D45' '36.66'
$ws.Range('E45').Value = '  -0.30%  '
# Row 46
Set-TextValue 'D46' '1.00'
# Row 47
Set-TextValue 'D47' '0.0264'
$ws.Range('E47').Value = '  +3.46%  '
# Row 48
Set-TextValue 'D48' '2.276.98'
$ws.Range('E48').Value = '  -0.32%  '
# Row 49
Set-TextValue 'D49' '1.00'
$ws.Range('E49').Value = '  +3.96%  '
# Row 50
Set-TextValue 'D50' '6.13'
$ws.Range('E50').Value = '  +2.45%  '
# Row 51
Set-TextValue 'D51' '20.71'
$ws.Range('E51').Value = '  -0.54%  '
